$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old third data row (duplicate "Calidad" escenario row) - it is no longer needed.
$ws.Rows("3:3").Delete()

# Insert a new column for "out_ref_externa" right before the old AN column,
# shifting out_n_pedido_sap / out_n_entrega / out_n_factura / out_n_orden_compra one column to the right.
$ws.Columns("AN:AN").Insert()
$ws.Columns("AN:AN").ColumnWidth = 22.592447916666668

$ws.Range("AN1").Value = "out_ref_externa"

# Update the remaining data row with the new "LOTE1125" escenario values.
$ws.Range("H2").Value = "SI"
$ws.Range("I2").Value = "LOTE1125"

# These values look numeric but must be stored as plain text (no quote-prefix marker),
# matching how the source cells were already typed as text in the workbook.
foreach ($pair in @(
        @{ Addr = "J2"; Val = "20200115" },
        @{ Addr = "K2"; Val = "20210906" },
        @{ Addr = "L2"; Val = "1125" },
        @{ Addr = "Y2"; Val = "20191116" },
        @{ Addr = "Z2"; Val = "3407077" }
    )) {
    $r = $ws.Range($pair.Addr)
    $r.NumberFormat = "@"
    $r.Value = $pair.Val
    $r.ClearFormats()
}

$ws.Range("N2").Value = "R0000-00001125"
$ws.Range("O2").Value = "OFF"
$ws.Range("M2").Value = 27
$ws.Range("AD2").Value = 27

$ws.Range("AI2").Value = "1002;LOTE1125;20210906;DESC LOTE1125;LOTE1125;N;LOTE1125;;AR;;;;;;N;;;210;20181103;;N;"
$ws.Range("AJ2").Value = "02000000001125;FNET;FNET;02;ZRET;20200115;;1800000122;20200115;08:00;16:00;20000;Remito electrónico Test;;;1002;27;C/U;LOTE1125;;;;;0000-00001125;1125;20200115;"
$ws.Range("AK2").Value = "PEDIDO            20200115002C001CLIENTESAPNROOC032202001151002              27           "
$ws.Range("AL2").Value = "202001150002073900PEDIDO                                                                          1002              27       139                         02        03      NUMEROWE "
$ws.Range("AM2").Value = "VTD02    133198CLIENTESAPC00120200115NROOC1  27     /////ESTE PEDIDO HA SIDO CREADO POR UN PROCESO DE AUTOMATIZACION./////                    1002              27     816 0  "

$ws.Range("AN2").Value = 3407087
# Leading apostrophe: force text with an explicit quote-prefix (matches original formatting).
$ws.Range("AO2").Value = "'0001128513"

$ws.Range("AN9").Select()
